$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 867.6667
$ws.Range("I31").Value = 867.6667
$ws.Range("K31").Value = 2603.0001
$ws.Range("M31").Value = -2373.0001
$ws.Range("H41").Value = 315.9524
$ws.Range("I41").Value = 272.1111
$ws.Range("J41").Value = 348.83334
$ws.Range("K41").Value = 272.1111
$ws.Range("L41").Value = 348.83334
$ws.Range("M41").Value = 167.8889
$ws.Range("N41").Value = -1228.83334
$ws.Range("H74").Value = 4113.476
$ws.Range("I74").Value = 3368.3
$ws.Range("J74").Value = 4790.909
$ws.Range("K74").Value = 3368.3
$ws.Range("L74").Value = 4790.909
$ws.Range("M74").Value = -2432.3
$ws.Range("N74").Value = -6662.909
$ws.Range("H77").Value = 4113.476
$ws.Range("I77").Value = 3368.3
$ws.Range("J77").Value = 4790.909
$ws.Range("K77").Value = 16841.5
$ws.Range("L77").Value = 23954.545
$ws.Range("M77").Value = -12161.5
$ws.Range("N77").Value = -33314.545
$ws.Range("H87").Value = 15626.313
$ws.Range("J87").Value = 15626.313
$ws.Range("L87").Value = 15626.313
$ws.Range("N87").Value = -18122.313
$ws.Range("H90").Value = 15626.313
$ws.Range("J90").Value = 15626.313
$ws.Range("L90").Value = 46878.939
$ws.Range("N90").Value = -59358.939
$ws.Range("H100").Value = 3318.3333
$ws.Range("I100").Value = 2997.9167
$ws.Range("J100").Value = 4600
$ws.Range("K100").Value = 2997.9167
$ws.Range("L100").Value = 4600
$ws.Range("M100").Value = -2456.9167
$ws.Range("N100").Value = -5682
$ws.Range("H103").Value = 581.5
$ws.Range("I103").Value = 551.3333
$ws.Range("K103").Value = 1653.9999
$ws.Range("M103").Value = -1067.9999
$ws.Range("H112").Value = 1226.8438
$ws.Range("J112").Value = 1273.5714
$ws.Range("L112").Value = 3820.7142
$ws.Range("N112").Value = -6036.7142
$ws.Range("H129").Value = 1023.5769
$ws.Range("J129").Value = 1037.8534
$ws.Range("L129").Value = 3113.5602
$ws.Range("N129").Value = -13113.5602
$ws.Range("H138").Value = 4630.324
$ws.Range("I138").Value = 2549.7856
$ws.Range("J138").Value = 5896.7393
$ws.Range("K138").Value = 7649.3568
$ws.Range("L138").Value = 17690.2179
$ws.Range("M138").Value = -2509.3568
$ws.Range("N138").Value = -27970.2179
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 800
$ws.Range("I41").Value = 800
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 800
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -386
$ws.Range("N41").ClearContents()
$ws.Range("H74").Value = 2263.6453
$ws.Range("I74").Value = 1497.8422
$ws.Range("J74").Value = 3476.1667
$ws.Range("K74").Value = 1497.8422
$ws.Range("L74").Value = 3476.1667
$ws.Range("M74").Value = -623.8422
$ws.Range("N74").Value = -5224.1667
$ws.Range("H77").Value = 2263.6453
$ws.Range("I77").Value = 1497.8422
$ws.Range("J77").Value = 3476.1667
$ws.Range("K77").Value = 7489.211
$ws.Range("L77").Value = 17380.8335
$ws.Range("M77").Value = -3121.211
$ws.Range("N77").Value = -26116.8335
$ws.Range("H97").Value = 673.8461
$ws.Range("I97").Value = 374.2857
$ws.Range("J97").Value = 1023.3333
$ws.Range("K97").Value = 374.2857
$ws.Range("L97").Value = 1023.3333
$ws.Range("M97").Value = 121.7143
$ws.Range("N97").Value = -2015.3333
$ws.Range("H102").Value = 2626.3333
$ws.Range("I102").Value = 1449.5
$ws.Range("K102").Value = 1449.5
$ws.Range("M102").Value = 172.5
$ws.Range("H132").Value = 2617.1135
$ws.Range("I132").Value = 2272.7917
$ws.Range("J132").Value = 3030.3
$ws.Range("K132").Value = 6818.375100000001
$ws.Range("L132").Value = 9090.900000000001
$ws.Range("M132").Value = -4288.375100000001
$ws.Range("N132").Value = -14150.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 528.1
$ws.Range("I94").Value = 538.8182
$ws.Range("K94").Value = 538.8182
$ws.Range("M94").Value = -87.81820000000005
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
$ws.Range("H140").Value = 58065
$ws.Range("J140").Value = 58065
$ws.Range("L140").Value = 58065
$ws.Range("N140").Value = -68425
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3910.4185
$ws.Range("I31").Value = 2102.0469
$ws.Range("J31").Value = 7314.4116
$ws.Range("K31").Value = 2102.0469
$ws.Range("L31").Value = 7314.4116
$ws.Range("M31").Value = -1807.0469
$ws.Range("N31").Value = -7904.4116
$ws.Range("H34").Value = 3910.4185
$ws.Range("I34").Value = 2102.0469
$ws.Range("J34").Value = 7314.4116
$ws.Range("K34").Value = 2102.0469
$ws.Range("L34").Value = 7314.4116
$ws.Range("M34").Value = -1900.0469
$ws.Range("N34").Value = -7718.4116
$ws.Range("H51").Value = 30747.477
$ws.Range("J51").Value = 31784.85
$ws.Range("L51").Value = 31784.85
$ws.Range("N51").Value = -33256.85
$ws.Range("H61").Value = 30747.477
$ws.Range("J61").Value = 31784.85
$ws.Range("L61").Value = 31784.85
$ws.Range("N61").Value = -32480.85
$ws.Range("H99").Value = 1461294.8
$ws.Range("I99").Value = 2135454.8
$ws.Range("J99").Value = 16666.285
$ws.Range("K99").Value = 2135454.8
$ws.Range("L99").Value = 16666.285
$ws.Range("M99").Value = -2133956.8
$ws.Range("N99").Value = -19662.285
$ws.Range("H122").Value = 1224058.6
$ws.Range("I122").Value = 167335.67
$ws.Range("K122").Value = 502007.01
$ws.Range("M122").Value = -499557.01
$ws.Range("H126").Value = 1461294.8
$ws.Range("I126").Value = 2135454.8
$ws.Range("J126").Value = 16666.285
$ws.Range("K126").Value = 6406364.399999999
$ws.Range("L126").Value = 49998.855
$ws.Range("M126").Value = -6403894.399999999
$ws.Range("N126").Value = -54938.855
$ws.Range("H140").Value = 78250.75
$ws.Range("J140").Value = 78250.75
$ws.Range("L140").Value = 78250.75
$ws.Range("N140").Value = -88610.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1232223
$ws.Range("I68").Value = 4000624.8
$ws.Range("J68").Value = 1822.2222
$ws.Range("K68").Value = 12001874.4
$ws.Range("L68").Value = 5466.6666
$ws.Range("M68").Value = -12001063.4
$ws.Range("N68").Value = -7088.6666
$ws.Range("H71").Value = 1232223
$ws.Range("I71").Value = 4000624.8
$ws.Range("J71").Value = 1822.2222
$ws.Range("K71").Value = 36005623.2
$ws.Range("L71").Value = 16399.9998
$ws.Range("M71").Value = -36001567.2
$ws.Range("N71").Value = -24511.9998
$ws.Range("H113").Value = 1263.5
$ws.Range("I113").Value = 1598.9
$ws.Range("J113").Value = 844.25
$ws.Range("K113").Value = 4796.700000000001
$ws.Range("L113").Value = 2532.75
$ws.Range("M113").Value = -2626.700000000001
$ws.Range("N113").Value = -6872.75
$ws.Range("H131").Value = 791.2
$ws.Range("I131").Value = 473
$ws.Range("J131").Value = 843
$ws.Range("K131").Value = 1419
$ws.Range("L131").Value = 2529
$ws.Range("M131").Value = 3621
$ws.Range("N131").Value = -12609
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1094.2106
$ws.Range("I2").Value = 952.7273
$ws.Range("K2").Value = 952.7273
$ws.Range("M2").Value = -839.7273
$ws.Range("H44").Value = 11999.5
$ws.Range("J44").Value = 11999.5
$ws.Range("L44").Value = 11999.5
$ws.Range("N44").Value = -13191.5
$ws.Range("H80").Value = 2166.6667
$ws.Range("I80").Value = 2111.875
$ws.Range("J80").Value = 2229.2856
$ws.Range("K80").Value = 2111.875
$ws.Range("L80").Value = 2229.2856
$ws.Range("M80").Value = -1113.875
$ws.Range("N80").Value = -4225.2856
$ws.Range("H83").Value = 2166.6667
$ws.Range("I83").Value = 2111.875
$ws.Range("J83").Value = 2229.2856
$ws.Range("K83").Value = 10559.375
$ws.Range("L83").Value = 11146.428
$ws.Range("M83").Value = -5567.375
$ws.Range("N83").Value = -21130.428
$ws.Range("H123").Value = 39746
$ws.Range("J123").Value = 39746
$ws.Range("L123").Value = 39746
$ws.Range("N123").Value = -44646
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 75020.42999999999
$ws.Range("I7").Value = 94116.91
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 94116.91
$ws.Range("L7").Value = 5000
$ws.Range("M7").Value = -94004.91
$ws.Range("N7").Value = -5224
$ws.Range("H22").Value = 956.25
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 930
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 930
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -1520
$ws.Range("H27").Value = 956.25
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 930
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 930
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -1144
$ws.Range("H46").Value = 1716.6
$ws.Range("I46").Value = 1924
$ws.Range("J46").Value = 1578.3334
$ws.Range("K46").Value = 1924
$ws.Range("L46").Value = 1578.3334
$ws.Range("M46").Value = -1736
$ws.Range("N46").Value = -1954.3334
$ws.Range("H80").Value = 19128
$ws.Range("J80").Value = 19128
$ws.Range("L80").Value = 19128
$ws.Range("N80").Value = -21374
$ws.Range("H83").Value = 19128
$ws.Range("J83").Value = 19128
$ws.Range("L83").Value = 57384
$ws.Range("N83").Value = -68616
$ws.Range("H93").Value = 1800
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -252
$ws.Range("N93").Value = -4496
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H122").Value = 2465.25
$ws.Range("I122").Value = 2250.5
$ws.Range("J122").Value = 2680
$ws.Range("K122").Value = 6751.5
$ws.Range("L122").Value = 8040
$ws.Range("M122").Value = -4301.5
$ws.Range("N122").Value = -12940
$ws.Range("H126").Value = 75020.42999999999
$ws.Range("I126").Value = 94116.91
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 282350.73
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -279880.73
$ws.Range("N126").Value = -19940
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 43333.332
$ws.Range("J51").Value = 43333.332
$ws.Range("L51").Value = 43333.332
$ws.Range("N51").Value = -44353.332
$ws.Range("H54").Value = 31544.445
$ws.Range("J54").Value = 31544.445
$ws.Range("L54").Value = 31544.445
$ws.Range("N54").Value = -32584.445
$ws.Range("H96").Value = 1986.9
$ws.Range("I96").Value = 2018.7778
$ws.Range("J96").Value = 1700
$ws.Range("K96").Value = 2018.7778
$ws.Range("L96").Value = 1700
$ws.Range("M96").Value = -645.7778000000001
$ws.Range("N96").Value = -4446
$ws.Range("H123").Value = 21111.111
$ws.Range("J123").Value = 21111.111
$ws.Range("L123").Value = 21111.111
$ws.Range("N123").Value = -30911.111
$ws.Range("H126").Value = 69378.52
$ws.Range("I126").Value = 72955.11
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 218865.33
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -216395.33
$ws.Range("N126").Value = -19940
$ws.Range("H132").Value = 1365.289
$ws.Range("I132").Value = 1123.5775
$ws.Range("J132").Value = 2268.5264
$ws.Range("K132").Value = 3370.7325
$ws.Range("L132").Value = 6805.5792
$ws.Range("M132").Value = -840.7325000000001
$ws.Range("N132").Value = -11865.5792
